$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.282.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.743.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +9.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.768.79"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.77"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.94%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.19%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.252.30"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.17"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.188.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.23%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.763.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.10"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.536"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.14"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.172"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.63"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0912"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +12.19%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.17%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +20.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "174.84"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.71"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.95%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.31%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +10.79%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.02%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "343.28"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.47"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.88"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.09"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0598"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.61%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.85"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0260"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.98%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.15%  "
